$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.844.22'
$ws.Cells.Item(2, 5).Value = '  +2.86%  '
$ws.Cells.Item(3, 4).Value = '2.606.95'
$ws.Cells.Item(3, 5).Value = '  +1.23%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '577.00'
$ws.Cells.Item(5, 5).Value = '  +4.06%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '143.66'
$ws.Cells.Item(6, 5).Value = '  +1.46%  '
$ws.Cells.Item(7, 5).Value = '  -0.30%  '
$ws.Cells.Item(8, 5).Value = '  +0.43%  '
$ws.Cells.Item(9, 4).Value = '2.632.75'
$ws.Cells.Item(9, 5).Value = '  +1.98%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.56'
$ws.Cells.Item(10, 5).Value = '  -2.22%  '
$ws.Cells.Item(11, 5).Value = '  +2.33%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.156'
$ws.Cells.Item(12, 5).Value = '  -4.94%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.370'
$ws.Cells.Item(13, 5).Value = '  +5.54%  '
$ws.Cells.Item(14, 4).Value = '3.069.94'
$ws.Cells.Item(14, 5).Value = '  +1.28%  '
$ws.Cells.Item(15, 4).Value = '60.798.89'
$ws.Cells.Item(15, 5).Value = '  +2.83%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '23.44'
$ws.Cells.Item(16, 5).Value = '  +2.06%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000143'
$ws.Cells.Item(17, 5).Value = '  +4.66%  '
$ws.Cells.Item(18, 4).Value = '2.619.95'
$ws.Cells.Item(18, 5).Value = '  +1.58%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.32'
$ws.Cells.Item(19, 5).Value = '  +9.74%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '4.67'
$ws.Cells.Item(20, 5).Value = '  +2.87%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '350.23'
$ws.Cells.Item(21, 5).Value = '  +3.99%  '
$ws.Cells.Item(22, 5).Value = '  +7.17%  '
$ws.Cells.Item(23, 5).Value = '  +0.01%  '
$ws.Cells.Item(24, 5).Value = '  +10.35%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '63.32'
$ws.Cells.Item(25, 5).Value = '  +1.18%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.997'
$ws.Cells.Item(26, 5).Value = '  -0.27%  '
$ws.Cells.Item(27, 5).Value = '  +1.58%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.93'
$ws.Cells.Item(28, 5).Value = '  +7.41%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0802'
$ws.Cells.Item(29, 5).Value = '  +4.25%  '
$ws.Cells.Item(30, 5).Value = '  +9.64%  '
$ws.Cells.Item(31, 2).Value = 'Aptos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '6.36'
$ws.Cells.Item(31, 5).Value = '  +2.75%  '
$ws.Cells.Item(32, 2).Value = 'USDe'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.997'
$ws.Cells.Item(32, 5).Value = '  -0.11%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '162.51'
$ws.Cells.Item(33, 5).Value = '  +2.39%  '
$ws.Cells.Item(34, 5).Value = '  +2.70%  '
$ws.Cells.Item(35, 5).Value = '  +16.61%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.28'
$ws.Cells.Item(36, 5).Value = '  +6.34%  '
$ws.Cells.Item(37, 5).Value = '  +6.73%  '
$ws.Cells.Item(38, 5).Value = '  +9.46%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '37.94'
$ws.Cells.Item(39, 5).Value = '  +1.81%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.89'
$ws.Cells.Item(40, 5).Value = '  +6.20%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.852'
$ws.Cells.Item(41, 5).Value = '  +0.18%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '302.63'
$ws.Cells.Item(42, 5).Value = '  +4.88%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '134.28'
$ws.Cells.Item(43, 5).Value = '  -1.97%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '19.99'
$ws.Cells.Item(44, 5).Value = '  +5.84%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.994'
$ws.Cells.Item(45, 5).Value = '  -0.48%  '
$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0986'
$ws.Cells.Item(46, 5).Value = '  +1.25%  '
$ws.Cells.Item(47, 2).Value = 'Mantle'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.608'
$ws.Cells.Item(47, 5).Value = '  +2.69%  '
$ws.Cells.Item(48, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '20.32'
$ws.Cells.Item(48, 5).Value = '  +9.27%  '
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '5.05'
$ws.Cells.Item(49, 5).Value = '  +12.28%  '
$ws.Cells.Item(50, 5).Value = '  +4.17%  '
$ws.Cells.Item(51, 5).Value = '  +4.16%  '
